# Insert a new weekly price record as row 70 (pushing the existing rows
# 70-86 down to 71-87), matching the "Fruta / hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70; Excel shifts rows 70-86 down to 71-87
# and copies formatting (incl. the date style on column D) from the row above.
$ws.Rows.Item(70).Insert()

$ws.Cells.Item(70, 1).Value = 1
$ws.Cells.Item(70, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(70, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(70, 4).Value = 44543
$ws.Cells.Item(70, 5).Value = 15
$ws.Cells.Item(70, 6).Value = 100112042
$ws.Cells.Item(70, 7).Value = "Locoto"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 150
$ws.Cells.Item(70, 11).Value = 14000
$ws.Cells.Item(70, 12).Value = 15000
$ws.Cells.Item(70, 13).Value = 14500
$ws.Cells.Item(70, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(70, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(70, 16).Value = 725
$ws.Cells.Item(70, 17).Value = 20
$ws.Cells.Item(70, 18).Value = "Hortaliza"
